$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13
$ws.Range("H21").Value = 8799.5
$ws.Range("I21").Value = 8799.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 8799.5
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -8331.5
$ws.Range("H23").Value = 8799.5
$ws.Range("I23").Value = 8799.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 8799.5
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -8565.5
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 100
$ws.Range("K29").Value = 300
$ws.Range("M29").Value = -19
$ws.Range("H58").Value = 2220.75
$ws.Range("J58").Value = 2294.6667
$ws.Range("L58").Value = 6884.000100000001
$ws.Range("N58").Value = -7184.000100000001
$ws.Range("H87").Value = 83902
$ws.Range("J87").Value = 83902
$ws.Range("L87").Value = 83902
$ws.Range("N87").Value = -86398
$ws.Range("H90").Value = 83902
$ws.Range("J90").Value = 83902
$ws.Range("L90").Value = 251706
$ws.Range("N90").Value = -264186
$ws.Range("H92").Value = 1054.9
$ws.Range("I92").Value = 972.1111
$ws.Range("J92").Value = 1800
$ws.Range("K92").Value = 972.1111
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = 275.8889
$ws.Range("N92").Value = -4296
$ws.Range("H98").Value = 4158
$ws.Range("I98").Value = 989.6
$ws.Range("K98").Value = 989.6
$ws.Range("M98").Value = 508.4
$ws.Range("H103").Value = 1000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 3000
$ws.Range("N103").Value = -4172
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -6746
$ws.Range("H122").Value = 4158
$ws.Range("I122").Value = 989.6
$ws.Range("K122").Value = 2968.8
$ws.Range("M122").Value = -518.8000000000002
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1230.2727
$ws.Range("I61").Value = 1282
$ws.Range("J61").Value = 997.5
$ws.Range("K61").Value = 1282
$ws.Range("L61").Value = 997.5
$ws.Range("M61").Value = -1070
$ws.Range("N61").Value = -1421.5
$ws.Range("H74").Value = 12467.556
$ws.Range("I74").Value = 14315.571
$ws.Range("J74").Value = 5999.5
$ws.Range("K74").Value = 14315.571
$ws.Range("L74").Value = 5999.5
$ws.Range("M74").Value = -13441.571
$ws.Range("N74").Value = -7747.5
$ws.Range("H77").Value = 12467.556
$ws.Range("I77").Value = 14315.571
$ws.Range("J77").Value = 5999.5
$ws.Range("K77").Value = 71577.855
$ws.Range("L77").Value = 29997.5
$ws.Range("M77").Value = -67209.855
$ws.Range("N77").Value = -38733.5
$ws.Range("H132").Value = 2868.4707
$ws.Range("I132").Value = 2272.25
$ws.Range("J132").Value = 4299.4
$ws.Range("K132").Value = 6816.75
$ws.Range("L132").Value = 12898.2
$ws.Range("M132").Value = -4286.75
$ws.Range("N132").Value = -17958.2
$ws.Range("H136").Value = 1230.2727
$ws.Range("I136").Value = 1282
$ws.Range("J136").Value = 997.5
$ws.Range("K136").Value = 3846
$ws.Range("L136").Value = 2992.5
$ws.Range("M136").Value = -1296
$ws.Range("N136").Value = -8092.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1638.6666
$ws.Range("I107").Value = 1693.5
$ws.Range("K107").Value = 1693.5
$ws.Range("M107").Value = 226.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 685.5714
$ws.Range("I22").Value = 839.8
$ws.Range("K22").Value = 839.8
$ws.Range("M22").Value = -489.8
$ws.Range("H74").Value = 68521.336
$ws.Range("J74").Value = 68521.336
$ws.Range("L74").Value = 68521.336
$ws.Range("N74").Value = -70269.336
$ws.Range("H77").Value = 68521.336
$ws.Range("J77").Value = 68521.336
$ws.Range("L77").Value = 205564.008
$ws.Range("N77").Value = -214300.008
$ws.Range("H99").Value = 7666.6665
$ws.Range("I99").Value = 7666.6665
$ws.Range("K99").Value = 7666.6665
$ws.Range("M99").Value = -6168.6665
$ws.Range("H126").Value = 7666.6665
$ws.Range("I126").Value = 7666.6665
$ws.Range("K126").Value = 22999.9995
$ws.Range("M126").Value = -20529.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -188
$ws.Range("H137").Value = 950
$ws.Range("I137").Value = 1155
$ws.Range("J137").Value = 540
$ws.Range("K137").Value = 3465
$ws.Range("L137").Value = 1620
$ws.Range("M137").Value = 1635
$ws.Range("N137").Value = -11820
$ws.Range("H138").Value = 2661.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H122").Value = 5500
$ws.Range("I122").Value = 5500
$ws.Range("K122").Value = 16500
$ws.Range("M122").Value = -14050
$ws.Range("H132").Value = 2243.4167
$ws.Range("I132").Value = 1992.3
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 5976.9
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -3446.9
$ws.Range("N132").Value = -15557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2984.6667
$ws.Range("I68").Value = 2984.6667
$ws.Range("K68").Value = 2984.6667
$ws.Range("M68").Value = -2235.6667
$ws.Range("H71").Value = 2984.6667
$ws.Range("I71").Value = 2984.6667
$ws.Range("K71").Value = 14923.3335
$ws.Range("M71").Value = -11179.3335
$ws.Range("H100").Value = 1911
$ws.Range("I100").Value = 1911
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1911
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -1370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7546.3335
$ws.Range("J62").Value = 7546.3335
$ws.Range("L62").Value = 7546.3335
$ws.Range("N62").Value = -8794.333500000001
$ws.Range("H65").Value = 7546.3335
$ws.Range("J65").Value = 7546.3335
$ws.Range("L65").Value = 37731.6675
$ws.Range("N65").Value = -43971.6675
$ws.Range("H126").Value = 2350
